{"js": "// Insert three new bullet paragraphs right after the \"Marketing Strategy and\n// Data-Driven Insights\" paragraph in the Siege Analytics / PARTNER section,\n// before the existing \"\u2022 Conducted comprehensive...\" bullet.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the exact target paragraph (there is only one occurrence in the doc).\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === \"Marketing Strategy and Data-Driven Insights\") {\n    target = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!target) {\n  throw new Error('Could not find paragraph \"Marketing Strategy and Data-Driven Insights\"');\n}\n\nconst newLines = [\n  \"\\u2022 Uncovered decades of demographic miscoding in voter files, discovering 500,000+ previously mischaracterized Democratic voters\",\n  \"\\u2022 Developed Python boundary estimation algorithm enabling mapping and analysis at every level of election in the United States\",\n  \"\\u2022 Algorithm reduced mapping costs by 75%, saving campaigns and organizations $5M+ and enabling smaller nonprofits to conduct redistricting analysis\"\n];\n\n// Insert each paragraph directly after the target, chaining so order is preserved.\nlet anchor = target;\nfor (const line of newLines) {\n  anchor = anchor.insertParagraph(line, Word.InsertLocation.after);\n}\n\nawait context.sync();\n", "ps1": "# Insert three new bullet paragraphs right after the \"Marketing Strategy and\n# Data-Driven Insights\" paragraph in the Siege Analytics / PARTNER section,\n# before the existing \"- Conducted comprehensive...\" bullet.\n\n$d = $word.ActiveDocument\n\n# Locate the 1-based paragraph index of the target paragraph (unique in doc).\n$targetIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text.TrimEnd([char]13)\n    if ($t -eq \"Marketing Strategy and Data-Driven Insights\") {\n        $targetIndex = $i\n        break\n    }\n}\nif ($targetIndex -eq -1) {\n    throw \"Could not find paragraph 'Marketing Strategy and Data-Driven Insights'\"\n}\n\n$lines = @(\n    \"\u2022 Uncovered decades of demographic miscoding in voter files, discovering 500,000+ previously mischaracterized Democratic voters\",\n    \"\u2022 Developed Python boundary estimation algorithm enabling mapping and analysis at every level of election in the United States\",\n    \"\u2022 Algorithm reduced mapping costs by 75%, saving campaigns and organizations `$5M+ and enabling smaller nonprofits to conduct redistricting analysis\"\n)\n\n$insertAfterIndex = $targetIndex\nforeach ($line in $lines) {\n    $p = $d.Paragraphs.Item($insertAfterIndex)\n    $p.Range.InsertParagraphAfter()\n    $newPara = $d.Paragraphs.Item($insertAfterIndex + 1)\n    $newPara.Range.Text = $line\n    $insertAfterIndex = $insertAfterIndex + 1\n}\n"}
